$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.697.86"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "2.122.51"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("D4").Value = "'1.012"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.92%  "
$ws.Range("D5").Value = "'338.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.22%  "
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("D7").Value = "'0.5257"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.68%  "
$ws.Range("D8").Value = "'0.4561"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.79%  "
$ws.Range("D9").Value = "'54.68"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.85%  "
$ws.Range("D10").Value = "'0.09130"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.28%  "
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("D12").Value = "'24.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "2.128.80"
$ws.Range("E13").Value = "  +1.73%  "
$ws.Range("D14").Value = "'6.863"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.52%  "
$ws.Range("D15").Value = "'8.163"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.91%  "
$ws.Range("D16").Value = "'0.00001179"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.00%  "
$ws.Range("D17").Value = "'97.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("D19").Value = "'0.06696"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").Value = "'19.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("D22").Value = "'6.325"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "30.766.84"
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("D24").Value = "'12.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.17%  "
$ws.Range("D25").Value = "'2.360"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("D26").Value = "2.368.16"
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("D27").Value = "'22.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.73%  "
$ws.Range("D28").Value = "'2.571"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "'164.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.83%  "
$ws.Range("D30").Value = "'134.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.11%  "
$ws.Range("E31").Value = "  +1.68%  "
$ws.Range("D32").Value = "'0.1076"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").Value = "'1.674"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("D34").Value = "'6.383"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.52%  "
$ws.Range("D35").Value = "'3.944"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.16%  "
$ws.Range("D36").Value = "'10.68"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.68%  "
$ws.Range("D37").Value = "'5.886"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.62%  "
$ws.Range("D38").Value = "'0.02644"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.02%  "
$ws.Range("D39").Value = "'0.06888"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.98%  "
$ws.Range("D40").Value = "'0.2332"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.06%  "
$ws.Range("D41").Value = "'12.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("D42").Value = "'0.6927"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("D43").Value = "'1.262"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("D44").Value = "'15.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.68%  "
$ws.Range("D45").Value = "'0.6505"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.44%  "
$ws.Range("E46").Value = "  +3.24%  "
$ws.Range("D47").Value = "'0.00000000370"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +21.85%  "
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("D49").Value = "'1.258"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("D50").Value = "'83.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.02%  "
$ws.Range("D51").Value = "'0.07305"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.66%  "

Write-Output "done"
